$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- G9 label changed from "无锡出院前一天" to "无锡一疗出院前一天" ---
$ws.Range("G9").Value = "无锡一疗出院前一天"

# --- Append a new data row (row 14), following the pattern of rows 4-13 ---

# B14:G14 get the same thin-border formatting used by the rest of the table
$ws.Range("B13:G13").Copy()
$ws.Range("B14:G14").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A14").Value = 45436
$ws.Range("A14").NumberFormat = "mm-dd-yy"  # built-in date format (numFmtId 14), no border
$ws.Range("B14").Value = 1.75
$ws.Range("C14").Value = 0.82
$ws.Range("D14").Value = 1.87
$ws.Range("E14").Value = 67
$ws.Range("F14").Value = 62
$ws.Range("G14").Value = "无锡二疗出院前一天"

# --- Update the view's selected cell (from G17 to F17) ---
$ws.Range("F17").Select()
